$wb = $excel.ActiveWorkbook

$oldTimestamp = "February 03 2026 17.29.55 EST"
$newTimestamp = "February 03 2026 18.05.36 EST"

$oldVersion = "Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on $oldTimestamp)"
$newVersion = "Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on $newTimestamp)"

# --- "About" sheet ---
$about = $wb.Worksheets.Item("About")

$about.Range("A2").Value = "Version: $newVersion"

$about.Range("A6").Value = "Recommended Citation:  `"Global Energy Monitor, Coal mine boundaries and methane sources for Cumberland Coal Mine, United States, M1012, version '$newVersion'. (See the CC license for attribution requirements if sharing or adapting the data set.)`""

# --- "Boundaries and methane sources" sheet ---
$data = $wb.Worksheets.Item("Boundaries and methane sources")

for ($row = 2; $row -le 31; $row++) {
    $cell = $data.Cells.Item($row, 19)  # column S = build_version
    $cell.Value = $newVersion
}
